# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Fri Nov 24 08:38:25 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.511.53'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '2.080.44'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''234.77'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").Value = '''0.626'
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '''57.54'
$ws.Range("E8").Value = '  -1.28%  '
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("D10").Value = '''0.0779'
$ws.Range("E10").Value = '  +2.32%  '
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").Value = '2.389.63'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '''14.43'
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = '''20.82'
$ws.Range("E14").Value = '  -1.41%  '
$ws.Range("D15").Value = '''0.782'
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("D16").Value = '''5.23'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '2.078.58'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '37.491.76'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").Value = '''6.21'
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").Value = '''69.70'
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("D21").Value = '0.0₃0821'
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("D22").Value = '''226.82'
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +1.72%  '
$ws.Range("D25").Value = '''2.40'
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("D26").Value = '''169.00'
$ws.Range("E26").Value = '  +1.35%  '
$ws.Range("D27").Value = '''8.91'
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  +3.87%  '
$ws.Range("E29").Value = '  -5.80%  '
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("D33").Value = '''0.0618'
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("D36").Value = '''3.37'
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("D39").Value = '''5.61'
$ws.Range("E39").Value = '  -4.29%  '
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.487.31'
$ws.Range("E41").Value = '  +2.11%  '
$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").Value = '''0.0952'
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("D43").Value = '''97.18'
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("E45").Value = '  -1.67%  '
$ws.Range("D46").Value = '''4.17'
$ws.Range("E46").Value = '  -9.88%  '
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("D48").Value = '''15.54'
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("D49").Value = '''7.26'
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("E50").Value = '  +1.40%  '
$ws.Range("D51").Value = '2.276.88'
$ws.Range("E51").Value = '  +0.73%  '
